# ContentIndexParser: Support for campaigns
#
# Adds a "my_campaign" worksheet + a new row/columns to the content_index
# sheet describing how to create a campaign.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) content_index: three new trailing columns (group / tags.1 / tags.2)
#    and a new summary row describing the campaign entry.
# ---------------------------------------------------------------------
$ci = $wb.Worksheets.Item("content_index")

$ci.Range("H1").Value = "group"
$ci.Range("I1").Value = "tags.1"
$ci.Range("J1").Value = "tags.2"

$ci.Range("I2").Value = "advanced"
$ci.Range("J2").Value = "type1"

$ci.Range("I3").Value = "advanced"
$ci.Range("J3").Value = "type2"

$ci.Range("I4").Value = "basic"

# new row describing the campaign-creation entry; gets its own style like
# the source workbook (distinct from the default row style)
$ci.Range("A6:I6").Style = "Normal"
$ci.Range("A6").Value = "create_campaign"
$ci.Range("B6").Value = "my_campaign"
$ci.Range("H6").Value = "My Group"
$ci.Range("I6").Value = "basic"

$ci.Range("H2").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) new sheet "my_campaign" appended at the end of the workbook
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$campaign = $wb.Worksheets.Add($null, $lastSheet)
$campaign.Name = "my_campaign"

$campaign.Range("A1").Value = "offset"
$campaign.Range("B1").Value = "unit"
$campaign.Range("C1").Value = "event_type"
$campaign.Range("D1").Value = "delivery_hour"
$campaign.Range("E1").Value = "message"
$campaign.Range("F1").Value = "relative_to"
$campaign.Range("G1").Value = "start_mode"
$campaign.Range("H1").Value = "flow"

$campaign.Range("A2").Value = 15
$campaign.Range("B2").Value = "H"
$campaign.Range("C2").Value = "F"
$campaign.Range("F2").Value = "Last Seen On"
$campaign.Range("G2").Value = "I"
$campaign.Range("H2").Value = "my_basic_flow"

$campaign.Range("A3").Value = 120
$campaign.Range("B3").Value = "D"
$campaign.Range("C3").Value = "M"
$campaign.Range("D3").Value = 12
$campaign.Range("E3").Value = "Hello!"
$campaign.Range("F3").Value = "Created On"
$campaign.Range("G3").Value = "S"

# leave the selection/active sheet on my_campaign, as in the target workbook
$campaign.Range("H47").Select() | Out-Null
